$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a boolean-typed "Unique" value (FALSE) to column L for rows 5, 6, 7
$range = $ws.Range("L5:L7")
$range.Value = $false
$range.NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# Move active selection to L7
$ws.Range("L7").Select()
